$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 443; this pushes the existing rows
# 443:463 down to 444:464 (and the sheet dimension grows to R464).
$ws.Rows("443:443").Insert()

# Populate the newly inserted row 443 with a fresh weekly price record.
# (Columns A, B, C, E, F, G, H, I, N, Q, R repeat the same constants
# used throughout this data block.)
$ws.Range("A443").Value = 10
$ws.Range("B443").Value = "Vega Modelo de Temuco"
$ws.Range("C443").Value = "La Araucanía"
$ws.Range("D443").Value = 45267
$ws.Range("E443").Value = 9
$ws.Range("F443").Value = 100112052
$ws.Range("G443").Value = "Albahaca"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Primera"
$ws.Range("J443").Value = 200
$ws.Range("K443").Value = 8000
$ws.Range("L443").Value = 9000
$ws.Range("M443").Value = 8750
$ws.Range("N443").Value = "$/paquete"
$ws.Range("O443").Value = "Región Metropolitana"
$ws.Range("P443").Value = 8750
$ws.Range("Q443").Value = 1
$ws.Range("R443").Value = "Hortaliza"

# Make sure the date cell keeps the same date number format as the
# rest of the column.
$ws.Range("D443").NumberFormat = $ws.Range("D444").NumberFormat
